$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "Image/audio/video/doc should be uploaded in the selected album and with the selected privacy only successfully."
$ws.Columns.Item(2).ColumnWidth = 75.0
$ws.Columns.Item(3).ColumnWidth = 71.83333333333333
$ws.Rows.Item(10).RowHeight = 45
$ws.Range("C10").Select()
